$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 10 (old "Nordeste" row) entirely - the data set shrinks from 9 data rows to 8
$ws.Rows.Item(10).Delete()

# Update the "Variável" column (B) for all data rows to the new period label
$newLabel = "Diferença 2024/10 - 2023/10"
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 2).Value = $newLabel
}

# Row 2: Amapá
$ws.Cells.Item(2, 1).Value = "Amapá"
$ws.Cells.Item(2, 3).Value = 5.49
$ws.Cells.Item(2, 4).Value = "1º"

# Row 3: Piauí
$ws.Cells.Item(3, 1).Value = "Piauí"
$ws.Cells.Item(3, 3).Value = 3.09
$ws.Cells.Item(3, 4).Value = "2º"

# Row 4: Sergipe
$ws.Cells.Item(4, 1).Value = "Sergipe"
$ws.Cells.Item(4, 3).Value = 2.77
$ws.Cells.Item(4, 4).Value = "3º"

# Row 5: Bahia
$ws.Cells.Item(5, 1).Value = "Bahia"
$ws.Cells.Item(5, 3).Value = 2.71
$ws.Cells.Item(5, 4).Value = "4º"

# Row 6: Ceará
$ws.Cells.Item(6, 1).Value = "Ceará"
$ws.Cells.Item(6, 3).Value = 2.26
$ws.Cells.Item(6, 4).Value = "5º"

# Row 7: Rio de Janeiro
$ws.Cells.Item(7, 1).Value = "Rio de Janeiro"
$ws.Cells.Item(7, 3).Value = 1.83
$ws.Cells.Item(7, 4).Value = "6º"

# Row 8: Brasil (no ranking)
$ws.Cells.Item(8, 1).Value = "Brasil"
$ws.Cells.Item(8, 3).Value = 1.24
$ws.Cells.Item(8, 4).Value = ""

# Row 9: Nordeste (no ranking)
$ws.Cells.Item(9, 1).Value = "Nordeste"
$ws.Cells.Item(9, 3).Value = 1.81
$ws.Cells.Item(9, 4).Value = ""
